# Automatische test-sync: 2025-08-03 02:21:50
# Adds a new log row (row 4) to the "Logs" sheet, mirroring the existing
# rows 2/3 test-mail entry, and bumps the matching count on "Dashboard".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$antwoord = "Beste klant,`nBedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, heb ik meer informatie nodig over wat u precies geregeld wilt hebben. Kunt u alstublieft wat meer details geven of specifiëren waar u hulp bij nodig heeft? Zodra ik meer informatie heb, zal ik ervoor zorgen dat dit zo spoedig mogelijk voor u wordt geregeld.`nMet vriendelijke groet,`n[E-mailassistent]"

$logs.Range("A4").Value = "Kun jij dit even regelen?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E4").Value = $antwoord
$logs.Range("F4").Value = "2025-08-03 02:21:09"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Nee"

# Dashboard: bump the "Intern verzoek / Actie voor medewerker" count (2 -> 3)
$dashboard.Range("B2").Value = 3

# Extend the conditional-formatting ranges so they also cover the new row 4
$ccols = @("D", "G", "H", "I", "J")
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + "3")
    $newRange = $logs.Range($col + "2:" + $col + "4")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
